$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ 'B' = 2; 'C' = 3; 'D' = 4; 'E' = 5 }

# Cells whose text is plain (non-numeric looking) - Coin names and Links
$bcUpdates = @(
    @{Row=6; Col='B'; Value='FTXToken'},
    @{Row=6; Col='C'; Value='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'},
    @{Row=7; Col='B'; Value='KuCoinToken'},
    @{Row=7; Col='C'; Value='https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'},
    @{Row=8; Col='B'; Value='MXToken'},
    @{Row=8; Col='C'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Row=9; Col='B'; Value='LiechtensteinCryptoassetsExchange'},
    @{Row=9; Col='C'; Value='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'},
    @{Row=10; Col='B'; Value='WazirX'},
    @{Row=10; Col='C'; Value='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'},
    @{Row=11; Col='B'; Value='MandalaExchangeToken'},
    @{Row=11; Col='C'; Value='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'},
    @{Row=12; Col='B'; Value='BitrueCoin'},
    @{Row=12; Col='C'; Value='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'},
    @{Row=13; Col='B'; Value='BitMartToken'},
    @{Row=13; Col='C'; Value='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Row=14; Col='B'; Value='BitForexToken'},
    @{Row=14; Col='C'; Value='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Row=15; Col='B'; Value='TigerCash'},
    @{Row=15; Col='C'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Row=16; Col='B'; Value='LEO'},
    @{Row=16; Col='C'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Row=17; Col='B'; Value='GateToken'},
    @{Row=17; Col='C'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'})

# Cells whose text looks numeric or percentage - Price and Volume columns.
# These must be forced to stay as text so the literal string (e.g. "0.1930")
# is preserved instead of being auto-converted into a number/percentage.
$deUpdates = @(
    @{Row=2; Col='D'; Value='311.95'},
    @{Row=2; Col='E'; Value='0.79%'},
    @{Row=3; Col='D'; Value='37.74'},
    @{Row=3; Col='E'; Value='-0.07%'},
    @{Row=4; Col='D'; Value='5.132'},
    @{Row=4; Col='E'; Value='0.51%'},
    @{Row=5; Col='D'; Value='0.07901'},
    @{Row=5; Col='E'; Value='0.59%'},
    @{Row=6; Col='D'; Value='1.904'},
    @{Row=6; Col='E'; Value='-2.73%'},
    @{Row=7; Col='D'; Value='8.276'},
    @{Row=7; Col='E'; Value='-0.28%'},
    @{Row=8; Col='D'; Value='0.9256'},
    @{Row=8; Col='E'; Value='-0.33%'},
    @{Row=9; Col='D'; Value='0.1203'},
    @{Row=9; Col='E'; Value='-10.93%'},
    @{Row=10; Col='D'; Value='0.1930'},
    @{Row=10; Col='E'; Value='-1.87%'},
    @{Row=11; Col='D'; Value='0.09111'},
    @{Row=11; Col='E'; Value='1.43%'},
    @{Row=12; Col='D'; Value='0.03325'},
    @{Row=12; Col='E'; Value='-4.24%'},
    @{Row=13; Col='D'; Value='0.09628'},
    @{Row=13; Col='E'; Value='-0.83%'},
    @{Row=14; Col='D'; Value='0.001383'},
    @{Row=14; Col='E'; Value='-0.52%'},
    @{Row=15; Col='D'; Value='0.005874'},
    @{Row=15; Col='E'; Value='-1.16%'},
    @{Row=16; Col='D'; Value='3.514'},
    @{Row=16; Col='E'; Value='-2.18%'},
    @{Row=17; Col='D'; Value='4.412'},
    @{Row=17; Col='E'; Value='0.96%'},
    @{Row=18; Col='D'; Value='3.099'},
    @{Row=18; Col='E'; Value='-0.14%'},
    @{Row=20; Col='D'; Value='5.288'},
    @{Row=20; Col='E'; Value='5.53%'},
    @{Row=21; Col='E'; Value='-1.73%'},
    @{Row=22; Col='E'; Value='2.96%'},
    @{Row=24; Col='D'; Value='0.04374'},
    @{Row=24; Col='E'; Value='0.64%'},
    @{Row=25; Col='D'; Value='0.001248'},
    @{Row=25; Col='E'; Value='2.23%'},
    @{Row=26; Col='D'; Value='0.004311'},
    @{Row=26; Col='E'; Value='-5.03%'},
    @{Row=27; Col='E'; Value='-9.75%'},
    @{Row=39; Col='D'; Value='0.02117'},
    @{Row=39; Col='E'; Value='-7.69%'},
    @{Row=40; Col='D'; Value='0.05183'},
    @{Row=40; Col='E'; Value='2.68%'},
    @{Row=41; Col='D'; Value='0.007667'},
    @{Row=41; Col='E'; Value='0.74%'},
    @{Row=42; Col='D'; Value='0.009048'},
    @{Row=42; Col='E'; Value='-8.25%'},
    @{Row=43; Col='D'; Value='0.1361'},
    @{Row=43; Col='E'; Value='0.44%'},
    @{Row=44; Col='D'; Value='0.002010'},
    @{Row=44; Col='E'; Value='-1.60%'},
    @{Row=45; Col='D'; Value='0.008605'},
    @{Row=45; Col='E'; Value='-2.16%'},
    @{Row=46; Col='D'; Value='0.00006699'},
    @{Row=46; Col='E'; Value='-1.98%'},
    @{Row=47; Col='E'; Value='-0.13%'},
    @{Row=48; Col='E'; Value='-7.87%'},
    @{Row=49; Col='D'; Value='0.002788'},
    @{Row=49; Col='E'; Value='-7.21%'},
    @{Row=50; Col='E'; Value='-0.13%'},
    @{Row=51; Col='E'; Value='-0.13%'})

foreach ($u in $bcUpdates) {
    $ws.Cells.Item($u.Row, $colIndex[$u.Col]).Value = $u.Value
}

# Force text format on the D/E cells we are about to touch so Excel does not
# reinterpret numeric-looking / percentage-looking strings as numbers.
$deRange = $ws.Range("D2:E51")
$deRange.NumberFormat = "@"

foreach ($u in $deUpdates) {
    $ws.Cells.Item($u.Row, $colIndex[$u.Col]).Value = $u.Value
}

# Restore the default (unstyled) look for those cells, now that the values
# are safely stored as text, so no stray number-format styling is left behind.
$deRange.Style = "Normal"
